$d = $word.ActiveDocument

# New "table of contents" style table: left column is the bold section
# heading, right column is the right-aligned page-range (e.g. "1:23").
$enDash = [char]0x2013
$entries = @(
    @{ Heading = "I INTRODUCTION AND OVERVIEW OF ALLEGATIONS"; Pages = "1:23" },
    @{ Heading = "II PARTIES"; Pages = "2:56" },
    @{ Heading = "III JURISDICTION AND VENUE"; Pages = "3:72" },
    @{ Heading = "IV FACTUAL ALLEGATIONS"; Pages = "3:81" },
    @{ Heading = "V SECTION 230 NON-IMMUNITY ALLEGATIONS"; Pages = "5:142" },
    @{ Heading = "VI CAUSE OF ACTION " + $enDash + " FRAUD"; Pages = "5:164" },
    @{ Heading = "VII PRAYER FOR RELIEF"; Pages = "6:181" }
)

function Escape-Xml([string]$text) {
    $text = $text.Replace("&", "&amp;")
    $text = $text.Replace("<", "&lt;")
    $text = $text.Replace(">", "&gt;")
    return $text
}

# Build the <w:tr> markup for every heading/page-range pair. Column widths
# (4320 twips each, half of a 8640-twip autofit table) and run formatting
# (bold + 24 half-points on the heading, non-bold + 24 half-points and
# right-justified on the page range) match the target table exactly.
$rowsXml = ""
foreach ($entry in $entries) {
    $heading = Escape-Xml $entry.Heading
    $pages = Escape-Xml $entry.Pages

    $rowsXml += "<w:tr>"
    $rowsXml += "<w:tc><w:tcPr><w:tcW w:type=""dxa"" w:w=""4320""/></w:tcPr>"
    $rowsXml += "<w:p><w:r><w:rPr><w:b/><w:sz w:val=""24""/></w:rPr><w:t>$heading</w:t></w:r></w:p></w:tc>"
    $rowsXml += "<w:tc><w:tcPr><w:tcW w:type=""dxa"" w:w=""4320""/></w:tcPr>"
    $rowsXml += "<w:p><w:pPr><w:jc w:val=""right""/></w:pPr><w:r><w:rPr><w:b w:val=""0""/><w:sz w:val=""24""/></w:rPr><w:t>$pages</w:t></w:r></w:p></w:tc>"
    $rowsXml += "</w:tr>"
}

$tableXml = "<w:tbl xmlns:w=""http://schemas.openxmlformats.org/wordprocessingml/2006/main"">"
$tableXml += "<w:tblPr><w:tblW w:type=""auto"" w:w=""0""/><w:tblLayout w:type=""autofit""/>"
$tableXml += "<w:tblLook w:firstColumn=""1"" w:firstRow=""1"" w:lastColumn=""0"" w:lastRow=""0"" w:noHBand=""0"" w:noVBand=""1"" w:val=""04A0""/></w:tblPr>"
$tableXml += "<w:tblGrid><w:gridCol w:w=""4320""/><w:gridCol w:w=""4320""/></w:tblGrid>"
$tableXml += $rowsXml
$tableXml += "</w:tbl>"

# Insert the table at the very end of the document body, right before the
# final section properties (mirrors where the diff places the new <w:tbl>:
# immediately after the existing "TABLE OF CONTENTS" paragraph).
$end = $d.Content
$end.Collapse(0)
$end.InsertXML($tableXml)

Write-Output ("Tables after insert: " + $d.Tables.Count)
Write-Output ("Rows in new table: " + $d.Tables(1).Rows.Count)
